# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a Price cell without Excel silently
# reinterpreting it as a number (which would strip meaningful trailing
# zeros / change the cell type). Temporarily forces a text format, then
# restores the cell's original style so no formatting changes leak in.
function Set-PriceText($addr, $text) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-PriceText 'D2' '25.931.83'
$ws.Range('E2').Value = '  -0.22%  '

Set-PriceText 'D3' '1.641.63'
$ws.Range('E3').Value = '  +0.35%  '

Set-PriceText 'D4' '1.005'
$ws.Range('E4').Value = '  -0.37%  '

Set-PriceText 'D5' '215.27'
$ws.Range('E5').Value = '  +0.03%  '

Set-PriceText 'D6' '0.5038'
$ws.Range('E6').Value = '  +0.00%  '

Set-PriceText 'D7' '1.004'
$ws.Range('E7').Value = '  -0.42%  '

Set-PriceText 'D8' '0.2571'
$ws.Range('E8').Value = '  +0.23%  '

Set-PriceText 'D9' '0.06416'
$ws.Range('E9').Value = '  +0.54%  '

Set-PriceText 'D10' '19.60'
$ws.Range('E10').Value = '  +0.64%  '

Set-PriceText 'D11' '0.07782'
$ws.Range('E11').Value = '  +0.62%  '

Set-PriceText 'D12' '4.270'
$ws.Range('E12').Value = '  +0.86%  '

Set-PriceText 'D13' '1.643.13'
$ws.Range('E13').Value = '  +0.21%  '

Set-PriceText 'D14' '1.868.60'
$ws.Range('E14').Value = '  +0.35%  '

Set-PriceText 'D15' '0.5429'
$ws.Range('E15').Value = '  +0.12%  '

Set-PriceText 'D16' '0.0₅7929'
$ws.Range('E16').Value = '  +0.22%  '

Set-PriceText 'D17' '64.47'
$ws.Range('E17').Value = '  +1.57%  '

Set-PriceText 'D18' '25.968.76'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('E19').Value = '  -0.33%  '

Set-PriceText 'D20' '199.19'
$ws.Range('E20').Value = '  -2.82%  '

Set-PriceText 'D21' '4.386'
$ws.Range('E21').Value = '  +1.30%  '

Set-PriceText 'D22' '9.908'
$ws.Range('E22').Value = '  -0.52%  '

Set-PriceText 'D23' '5.968'
$ws.Range('E23').Value = '  +0.42%  '

$ws.Range('E24').Value = '  -0.33%  '

Set-PriceText 'D25' '1.870'
$ws.Range('E25').Value = '  -5.08%  '

Set-PriceText 'D26' '141.14'
$ws.Range('E26').Value = '  -0.57%  '

Set-PriceText 'D27' '0.1134'
$ws.Range('E27').Value = '  -1.71%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-PriceText 'D28' '6.812'
$ws.Range('E28').Value = '  -0.19%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-PriceText 'D29' '15.70'
$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('E30').Value = '  +0.40%  '

Set-PriceText 'D31' '0.04928'
$ws.Range('E31').Value = '  -1.34%  '

Set-PriceText 'D32' '3.265'
$ws.Range('E32').Value = '  +0.23%  '

Set-PriceText 'D33' '3.205'
$ws.Range('E33').Value = '  +0.65%  '

Set-PriceText 'D34' '1.540'
$ws.Range('E34').Value = '  +0.33%  '

Set-PriceText 'D35' '2.370'
$ws.Range('E35').Value = '  +1.56%  '

Set-PriceText 'D36' '0.8929'
$ws.Range('E36').Value = '  -0.99%  '

$ws.Range('E37').Value = '  -1.17%  '

Set-PriceText 'D38' '1.149.29'
$ws.Range('E38').Value = '  +2.27%  '

Set-PriceText 'D39' '0.5555'
$ws.Range('E39').Value = '  -1.62%  '

Set-PriceText 'D40' '0.01570'
$ws.Range('E40').Value = '  +0.94%  '

$ws.Range('E41').Value = '  -0.34%  '

Set-PriceText 'D42' '5.720'
$ws.Range('E42').Value = '  +1.66%  '

Set-PriceText 'D43' '0.8115'
$ws.Range('E43').Value = '  +0.07%  '

Set-PriceText 'D44' '99.87'
$ws.Range('E44').Value = '  +0.18%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-PriceText 'D45' '1.779.55'
$ws.Range('E45').Value = '  +0.28%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-PriceText 'D46' '0.0₈119'
$ws.Range('E46').Value = '  +4.56%  '

Set-PriceText 'D47' '0.4522'
$ws.Range('E47').Value = '  -0.20%  '

Set-PriceText 'D48' '1.004'
$ws.Range('E48').Value = '  -0.19%  '

Set-PriceText 'D49' '54.72'
$ws.Range('E49').Value = '  -0.04%  '

Set-PriceText 'D50' '0.05041'
$ws.Range('E50').Value = '  -0.15%  '

Set-PriceText 'D51' '1.005'
$ws.Range('E51').Value = '  -0.08%  '
